$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Style = $ws.Range("F1").Style

$ws.Range("A6").Value = "tarun"
$ws.Range("B6").Value = "U654341"
$ws.Range("C6").Value = "2200039159@kluniversity.in"
$ws.Range("D6").Value = "Category-2"
$ws.Range("E6").Value = 2020

$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:2200039159@kluniversity.in")
$ws.Range("C6").Style = $ws.Range("C5").Style

$ws.Range("E6").Select()
